$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "0001_slr0611_right"
$ws.Range("B12").Value = "Wrap-Up-doc.docx"

$ws.Range("A13").Value = "sll0199_right"
$ws.Range("B13").Value = "flanks_short.xlsx"

$ws.Range("A13").Select()
